$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the next bitacora entry (row 31, merged C31:C32 / D31:D32 / E31:E32)
$ws.Range("C31").Value = "Controller"

# Give D31 the same date formatting as the other date cells in the column,
# then write the date value
$ws.Range("D7").Copy()
$ws.Range("D31").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("D31").Value = 44713

$ws.Range("E31").Value = "Creacion de controller de administrador"

# Move the active selection to the newly filled row, like the user did
# right after typing the new entry
$ws.Range("E31:E32").Select()
